# Models_v1.xlsx update: prune the model list down to the models that are
# actually still in use, removing the rows for models that were dropped
# from the analysis (and the accidental duplicate rows at the bottom of
# the sheet).
#
# Rows (1-based, as they exist in the sheet before any deletions) that are
# being removed, identified by the model name in column A:
#   6  CNRM-CM6-1
#   7  EC-Earth3-AerChem
#   9  EC-Earth3-CC
#   10 EC-Earth3-Veg-LR
#   11 EC-Earth3-Veg
#   12 HadGEM3-GC31-LL
#   15 MIROC-ES2H
#   20 NorESM2-MM   (duplicate trailing row)
#
# Deleting from the bottom of the sheet upward keeps the row numbers for
# rows we haven't processed yet stable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("20").Delete()
$ws.Rows("15").Delete()
$ws.Rows("12").Delete()
$ws.Rows("11").Delete()
$ws.Rows("10").Delete()
$ws.Rows("9").Delete()
$ws.Rows("7").Delete()
$ws.Rows("6").Delete()

# B3 (BCC-CSM2-MR's Ref date) had picked up a stray "apply number format"
# style along the way; restore it to the sheet's default (unstyled) look.
$ws.Range("B3").Style = "Normal"

# Leave the selection where it landed after the last deletion (the now
# -empty row below the shrunken table), matching the saved selection.
$ws.Rows("13").Select() | Out-Null
